# Gate_Planning.xlsx base model update
# - Flight Schedule sheet: replace placeholder flight data (A1/A2/A3, generic
#   pax counts) with real flight numbers, pax counts, and ETA/ETD headers.
# - Gates sheet: content unchanged, just move the stale selection.

$wb = $excel.ActiveWorkbook

# --- Flight Schedule sheet --------------------------------------------------
$flights = $wb.Worksheets.Item("Flight Schedule")

$flights.Range("A1").Value = "Flight No."
$flights.Range("B1").Value = "Pax"
$flights.Range("C1").Value = "ETA"
$flights.Range("D1").Value = "ETD"

$flights.Range("A2").Value = "KL2020"
$flights.Range("B2").Value = 110

$flights.Range("A3").Value = "KL358"
$flights.Range("B3").Value = 436

$flights.Range("A4").Value = "HV2587"
$flights.Range("B4").Value = 186

# --- Gates sheet -------------------------------------------------------------
$gates = $wb.Worksheets.Item("Gates")
$gates.Range("A2").Select()

# Leave the active sheet/selection on Flight Schedule!C2, matching the saved
# file's last selection.
$flights.Activate()
$flights.Range("C2").Select()
